{"js": "// The document contains an embedded SQL query (as plain text) inside a\n// table. Two lines of that query were edited:\n//\n// 1) The \"CityCode\" projection was simplified from a DECODE() call down\n//    to a straight column reference:\n//      DECODE(s1.\"CityCode\", '85', '96', S1.\"CityCode\") AS \"CityCode\",\n//      -> s1.\"CityCode\" AS \"CityCode\",\n//\n// 2) The denominator of the overdue ratio calculation gained an extra\n//    SUM(s1.\"OvduBal\") term:\n//      round((SUM(s1.\"OvduBal\") + SUM(s1.\"ColBal\")) /(SUM(s1.\"LoanBal\") + SUM(s1.\"ColBal\")), 4)\n//      -> round((SUM(s1.\"OvduBal\") + SUM(s1.\"ColBal\")) /(SUM(s1.\"LoanBal\") + SUM(s1.\"OvduBal\") + SUM(s1.\"ColBal\")), 4)\n\nconst body = context.document.body;\n\n// --- Edit 1: replace the DECODE(...) expression with a bare column ref.\nconst decodeResults = body.search(\n  'DECODE(s1.\"CityCode\", \\'85\\', \\'96\\', S1.\"CityCode\")',\n  { matchCase: true }\n);\ndecodeResults.load(\"text\");\nawait context.sync();\n\nif (decodeResults.items.length > 0) {\n  decodeResults.items[0].insertText('s1.\"CityCode\"', Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Edit 2: insert the extra SUM(s1.\"OvduBal\") term into the round() denominator.\nconst roundResults = body.search(\n  'round((SUM(s1.\"OvduBal\") + SUM(s1.\"ColBal\")) /(SUM(s1.\"LoanBal\") + SUM(s1.\"ColBal\")), 4)',\n  { matchCase: true }\n);\nroundResults.load(\"text\");\nawait context.sync();\n\nif (roundResults.items.length > 0) {\n  roundResults.items[0].insertText(\n    'round((SUM(s1.\"OvduBal\") + SUM(s1.\"ColBal\")) /(SUM(s1.\"LoanBal\") + SUM(s1.\"OvduBal\") + SUM(s1.\"ColBal\")), 4)',\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n", "ps1": "# The document contains an embedded SQL query (as plain text) inside a\n# table. Two lines of that query were edited:\n#\n# 1) The \"CityCode\" projection was simplified from a DECODE() call down\n#    to a straight column reference:\n#      DECODE(s1.\"CityCode\", '85', '96', S1.\"CityCode\") AS \"CityCode\",\n#      -> s1.\"CityCode\" AS \"CityCode\",\n#\n# 2) The denominator of the overdue ratio calculation gained an extra\n#    SUM(s1.\"OvduBal\") term:\n#      round((SUM(s1.\"OvduBal\") + SUM(s1.\"ColBal\")) /(SUM(s1.\"LoanBal\") + SUM(s1.\"ColBal\")), 4)\n#      -> round((SUM(s1.\"OvduBal\") + SUM(s1.\"ColBal\")) /(SUM(s1.\"LoanBal\") + SUM(s1.\"OvduBal\") + SUM(s1.\"ColBal\")), 4)\n#\n# NOTE: we assign Range.Text directly (instead of using Find.Replacement /\n# Find.Execute's Replace argument) because routing the replacement string\n# through Find.Replacement.Text triggers Word's \"smart quotes\" AutoCorrect\n# and turns the straight double quotes in the SQL text into curly ones.\n\n$d = $word.ActiveDocument\n\n# --- Edit 1: replace the DECODE(...) expression with a bare column ref.\n$range1 = $d.Content\n$find1 = $range1.Find\n$find1.ClearFormatting()\n$find1.Text = \"DECODE(s1.`\"CityCode`\", '85', '96', S1.`\"CityCode`\")\"\n$found1 = $find1.Execute()\nif ($found1) {\n    $range1.Text = \"s1.`\"CityCode`\"\"\n}\n\n# --- Edit 2: insert the extra SUM(s1.\"OvduBal\") term into the round() denominator.\n$range2 = $d.Content\n$find2 = $range2.Find\n$find2.ClearFormatting()\n$find2.Text = \"round((SUM(s1.`\"OvduBal`\") + SUM(s1.`\"ColBal`\")) /(SUM(s1.`\"LoanBal`\") + SUM(s1.`\"ColBal`\")), 4)\"\n$found2 = $find2.Execute()\nif ($found2) {\n    $range2.Text = \"round((SUM(s1.`\"OvduBal`\") + SUM(s1.`\"ColBal`\")) /(SUM(s1.`\"LoanBal`\") + SUM(s1.`\"OvduBal`\") + SUM(s1.`\"ColBal`\")), 4)\"\n}\n"}
